$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the angle-of-rotation matrix values that changed between revisions
# (each row below lists only the cells whose value differs from the original).

# Row 1
$ws.Range("B1").Value = 180
$ws.Range("C1").Value = 315
$ws.Range("D1").Value = 0
$ws.Range("G1").Value = 45
$ws.Range("K1").Value = 270
$ws.Range("M1").Value = 0
$ws.Range("N1").Value = 225
$ws.Range("O1").Value = 135
$ws.Range("P1").Value = 225
$ws.Range("Q1").Value = 135
$ws.Range("S1").Value = 135
$ws.Range("T1").Value = 270
$ws.Range("U1").Value = 45
$ws.Range("V1").Value = 270
$ws.Range("W1").Value = 180
$ws.Range("X1").Value = 135
$ws.Range("Y1").Value = 225
$ws.Range("Z1").Value = 270
$ws.Range("AC1").Value = 225

# Row 2
$ws.Range("C2").Value = 180
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 315
$ws.Range("G2").Value = 180
$ws.Range("H2").Value = 270
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 135
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 315
$ws.Range("N2").Value = 0
$ws.Range("P2").Value = 315
$ws.Range("Q2").Value = 180
$ws.Range("S2").Value = 180
$ws.Range("T2").Value = 315
$ws.Range("V2").Value = 315
$ws.Range("W2").Value = 135
$ws.Range("AA2").Value = 135
$ws.Range("AC2").Value = 0

# Row 3
$ws.Range("D3").Value = 225
$ws.Range("F3").Value = 225
$ws.Range("G3").Value = 225
$ws.Range("H3").Value = 45
$ws.Range("I3").Value = 45
$ws.Range("J3").Value = 45
$ws.Range("M3").Value = 225
$ws.Range("N3").Value = 45
$ws.Range("O3").Value = 225
$ws.Range("P3").Value = 45
$ws.Range("Q3").Value = 225
$ws.Range("R3").Value = 225
$ws.Range("S3").Value = 90
$ws.Range("T3").Value = 225
$ws.Range("U3").Value = 225
$ws.Range("V3").Value = 225
$ws.Range("W3").Value = 45
$ws.Range("X3").Value = 225
$ws.Range("Y3").Value = 225
$ws.Range("Z3").Value = 45
$ws.Range("AA3").Value = 45
$ws.Range("AB3").Value = 45
$ws.Range("AC3").Value = 45

# Row 4
$ws.Range("E4").Value = 45
$ws.Range("F4").Value = 315
$ws.Range("G4").Value = 45
$ws.Range("H4").Value = 135
$ws.Range("I4").Value = 135
$ws.Range("J4").Value = 225
$ws.Range("K4").Value = 135
$ws.Range("L4").Value = 45
$ws.Range("M4").Value = 315
$ws.Range("N4").Value = 135
$ws.Range("O4").Value = 315
$ws.Range("P4").Value = 315
$ws.Range("Q4").Value = 315
$ws.Range("R4").Value = 315
$ws.Range("T4").Value = 315
$ws.Range("V4").Value = 315
$ws.Range("W4").Value = 135
$ws.Range("X4").Value = 315
$ws.Range("Y4").Value = 135
$ws.Range("Z4").Value = 135
$ws.Range("AB4").Value = 135
$ws.Range("AC4").Value = 135

# Row 5
$ws.Range("G5").Value = 180
$ws.Range("I5").Value = 45
$ws.Range("P5").Value = 225
$ws.Range("U5").Value = 225
$ws.Range("V5").Value = 225
$ws.Range("Y5").Value = 45
$ws.Range("AA5").Value = 45

# Row 6
$ws.Range("H6").Value = 225
$ws.Range("I6").Value = 315
$ws.Range("K6").Value = 90
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 180
$ws.Range("P6").Value = 270
$ws.Range("Q6").Value = 135
$ws.Range("R6").Value = 0
$ws.Range("T6").Value = 270
$ws.Range("U6").Value = 270
$ws.Range("V6").Value = 270
$ws.Range("W6").Value = 90
$ws.Range("X6").Value = 0
$ws.Range("Y6").Value = 90
$ws.Range("AA6").Value = 90
$ws.Range("AB6").Value = 180

# Row 7
$ws.Range("H7").Value = 225
$ws.Range("J7").Value = 180
$ws.Range("K7").Value = 315
$ws.Range("L7").Value = 180
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 225
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 225
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 135
$ws.Range("U7").Value = 45
$ws.Range("V7").Value = 135
$ws.Range("W7").Value = 225
$ws.Range("X7").Value = 135
$ws.Range("Y7").Value = 45
$ws.Range("Z7").Value = 225
$ws.Range("AA7").Value = 225
$ws.Range("AB7").Value = 225
$ws.Range("AC7").Value = 225

# Row 8
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 180
$ws.Range("N8").Value = 225
$ws.Range("O8").Value = 180
$ws.Range("Q8").Value = 180
$ws.Range("S8").Value = 225
$ws.Range("T8").Value = 45
$ws.Range("U8").Value = 180
$ws.Range("V8").Value = 45
$ws.Range("W8").Value = 90
$ws.Range("Y8").Value = 225
$ws.Range("Z8").Value = 0
$ws.Range("AB8").Value = 225
$ws.Range("AC8").Value = 0

# Row 9
$ws.Range("J9").Value = 45
$ws.Range("K9").Value = 270
$ws.Range("M9").Value = 45
$ws.Range("N9").Value = 315
$ws.Range("O9").Value = 45
$ws.Range("Q9").Value = 45
$ws.Range("R9").Value = 90
$ws.Range("S9").Value = 45
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 90
$ws.Range("W9").Value = 180
$ws.Range("X9").Value = 90
$ws.Range("Y9").Value = 270
$ws.Range("Z9").Value = 270
$ws.Range("AA9").Value = 180
$ws.Range("AB9").Value = 270
$ws.Range("AC9").Value = 270

# Row 10
$ws.Range("K10").Value = 180
$ws.Range("L10").Value = 0
$ws.Range("P10").Value = 315
$ws.Range("R10").Value = 315
$ws.Range("T10").Value = 0
$ws.Range("U10").Value = 180
$ws.Range("V10").Value = 180
$ws.Range("W10").Value = 135
$ws.Range("X10").Value = 315
$ws.Range("Z10").Value = 135
$ws.Range("AA10").Value = 135
$ws.Range("AB10").Value = 0
$ws.Range("AC10").Value = 135

# Row 11
$ws.Range("L11").Value = 45
$ws.Range("M11").Value = 225
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("Q11").Value = 225
$ws.Range("R11").Value = 135
$ws.Range("S11").Value = 225
$ws.Range("W11").Value = 270
$ws.Range("Y11").Value = 180
$ws.Range("Z11").Value = 0
$ws.Range("AA11").Value = 45
$ws.Range("AC11").Value = 0

# Row 12
$ws.Range("M12").Value = 135
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 180
$ws.Range("P12").Value = 180
$ws.Range("R12").Value = 135
$ws.Range("S12").Value = 270
$ws.Range("U12").Value = 135
$ws.Range("V12").Value = 315
$ws.Range("W12").Value = 315
$ws.Range("X12").Value = 180
$ws.Range("Y12").Value = 180
$ws.Range("Z12").Value = 0
$ws.Range("AA12").Value = 315
$ws.Range("AB12").Value = 315
$ws.Range("AC12").Value = 0

# Row 13
$ws.Range("N13").Value = 135
$ws.Range("O13").Value = 315
$ws.Range("P13").Value = 270
$ws.Range("Q13").Value = 135
$ws.Range("S13").Value = 0
$ws.Range("T13").Value = 90
$ws.Range("V13").Value = 90
$ws.Range("Y13").Value = 90
$ws.Range("AC13").Value = 135

# Row 14
$ws.Range("O14").Value = 180
$ws.Range("P14").Value = 135
$ws.Range("R14").Value = 180
$ws.Range("S14").Value = 270
$ws.Range("T14").Value = 135
$ws.Range("U14").Value = 135
$ws.Range("V14").Value = 135
$ws.Range("W14").Value = 135
$ws.Range("X14").Value = 315
$ws.Range("Y14").Value = 315
$ws.Range("AA14").Value = 0

# Row 15
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = 90
$ws.Range("S15").Value = 90
$ws.Range("T15").Value = 180
$ws.Range("U15").Value = 180
$ws.Range("V15").Value = 180
$ws.Range("W15").Value = 0
$ws.Range("X15").Value = 90
$ws.Range("Y15").Value = 45
$ws.Range("AB15").Value = 0

# Row 16
$ws.Range("Q16").Value = 180
$ws.Range("T16").Value = 180
$ws.Range("U16").Value = 180
$ws.Range("V16").Value = 180
$ws.Range("X16").Value = 180
$ws.Range("Y16").Value = 270
$ws.Range("Z16").Value = 45
$ws.Range("AA16").Value = 45
$ws.Range("AB16").Value = 0
$ws.Range("AC16").Value = 45

# Row 17
$ws.Range("R17").Value = 45
$ws.Range("S17").Value = 270
$ws.Range("U17").Value = 135
$ws.Range("X17").Value = 45
$ws.Range("Z17").Value = 315
$ws.Range("AC17").Value = 90

# Row 18
$ws.Range("T18").Value = 90
$ws.Range("U18").Value = 0
$ws.Range("V18").Value = 90
$ws.Range("X18").Value = 90
$ws.Range("Z18").Value = 270
$ws.Range("AA18").Value = 315
$ws.Range("AB18").Value = 270

# Row 19
$ws.Range("T19").Value = 270
$ws.Range("X19").Value = 0
$ws.Range("AA19").Value = 180
$ws.Range("AB19").Value = 180

# Row 20
$ws.Range("U20").Value = 0
$ws.Range("V20").Value = 0
$ws.Range("W20").Value = 180
$ws.Range("X20").Value = 225
$ws.Range("Y20").Value = 225
$ws.Range("Z20").Value = 180
$ws.Range("AA20").Value = 180
$ws.Range("AC20").Value = 180

# Row 21
$ws.Range("V21").Value = 225
$ws.Range("Y21").Value = 270
$ws.Range("Z21").Value = 180
$ws.Range("AB21").Value = 270

# Row 22
$ws.Range("W22").Value = 270
$ws.Range("AA22").Value = 270

# Row 23
$ws.Range("Z23").Value = 45

# Row 24
$ws.Range("Y24").Value = 135
$ws.Range("Z24").Value = 90
$ws.Range("AB24").Value = 180

# Row 25
$ws.Range("AB25").Value = 315

# Row 26
$ws.Range("AB26").Value = 45
$ws.Range("AC26").Value = 45

# Row 27
$ws.Range("AB27").Value = 135

# Row 28
$ws.Range("AC28").Value = 315

# Update the selection / zoom to match the saved view state.
$ws.Range("A21:AC29").Select()
$excel.ActiveWindow.Zoom = 70
